# The "年号" (era-name) sheet lists Chinese dynasty/era/ruler rows.
# Two corrections to the Later Zhao (后赵) block around row 93:
#   1. A missing ruler, 石世 (Shi Shi), is inserted as a new row before
#      the existing row 93 (石遵), shifting everything below down by one.
#   2. The row that used to be row 94 (originally mislabeled 石祗 with
#      years 350-350) becomes row 95 after the insert, and is corrected
#      to be 石鉴 with years 349-350. The row that is now 94 (old row 93,
#      石遵) also gets its end-year corrected from 350 to 349.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 93 (pushes old rows 93-108 down to 94-109)
$ws.Rows("93:93").Insert()

# New row 93: 后赵 / 太宁 / 石世 / 349 / 349
$ws.Range("A93").Value = "后赵"
$ws.Range("B93").Value = "太宁"
$ws.Range("C93").Value = "石世"
$ws.Range("D93").Value = 349
$ws.Range("E93").Value = 349

# Row 94 (old row 93, 石遵): fix end year 350 -> 349
$ws.Range("E94").Value = 349

# Row 95 (old row 94): fix name 石祗 -> 石鉴 and start year 350 -> 349
$ws.Range("C95").Value = "石鉴"
$ws.Range("D95").Value = 349

# Update the sheet's active selection to match
$ws.Range("G95").Select()
